# Update "想去人数" (attendee count) figures on the "展览" and "全部类型"
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# (SheetName, Row, NewValue)
$updates = @(
    @{ Sheet = "展览";     Row = 3;  Value = 1342 },
    @{ Sheet = "展览";     Row = 5;  Value = 100 },
    @{ Sheet = "展览";     Row = 6;  Value = 59 },
    @{ Sheet = "展览";     Row = 8;  Value = 11598 },
    @{ Sheet = "展览";     Row = 12; Value = 23 },
    @{ Sheet = "展览";     Row = 16; Value = 141 },
    @{ Sheet = "展览";     Row = 18; Value = 5051 },
    @{ Sheet = "展览";     Row = 21; Value = 11325 },
    @{ Sheet = "展览";     Row = 28; Value = 18 },

    @{ Sheet = "全部类型"; Row = 3;  Value = 1342 },
    @{ Sheet = "全部类型"; Row = 5;  Value = 100 },
    @{ Sheet = "全部类型"; Row = 6;  Value = 59 },
    @{ Sheet = "全部类型"; Row = 8;  Value = 11598 },
    @{ Sheet = "全部类型"; Row = 12; Value = 23 },
    @{ Sheet = "全部类型"; Row = 17; Value = 141 },
    @{ Sheet = "全部类型"; Row = 19; Value = 5051 },
    @{ Sheet = "全部类型"; Row = 22; Value = 11325 },
    @{ Sheet = "全部类型"; Row = 29; Value = 18 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Cells.Item($u.Row, 6).Value = $u.Value
}
